# Generate Report for Handoff
# The b.md entry (row 3 on each sheet) is now "Ready for handoff" with a
# fresh handoff artifact (b.63290e5768f688058c7b37413b0a5c26c308f864.*)
# and an updated handoff datetime, on both the Overview sheet and the
# per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay($ws, $row, $col, $text) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Row -eq $row -And $h.Range.Column -eq $col) {
            $h.TextToDisplay = $text
        }
    }
}

# --- Overview sheet: row 3 is the "b.md" file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-26-12 14:26:35"

# --- zh-cn sheet: row 3 is the "b.md" source file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-12 14:26:32"
Set-HyperlinkDisplay $wsZhCn 3 4 "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

# --- de-de sheet: row 3 is the "b.md" source file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-12 14:26:35"
Set-HyperlinkDisplay $wsDeDe 3 4 "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
